$d = $word.ActiveDocument

function Add-EmptyPara($doc) {
    $cur = $doc.Paragraphs.Last.Range
    $cur.InsertParagraphAfter()
}

function Add-TextPara($doc, $text) {
    $cur = $doc.Paragraphs.Last.Range
    $cur.InsertParagraphAfter()
    $newPara = $doc.Paragraphs.Last
    $newPara.Range.InsertAfter($text)
}

function Merge-WithPrevious($doc) {
    # Deletes the paragraph mark that ends the paragraph immediately before
    # the last paragraph, joining it with the last paragraph while keeping
    # their runs distinct (mirrors Word's "delete pilcrow to merge" editing
    # behaviour, which does not coalesce differently-authored runs).
    $count = $doc.Paragraphs.Count
    $prevPara = $doc.Paragraphs.Item($count - 1)
    $markStart = $prevPara.Range.End - 1
    $markRange = $doc.Range($markStart, $markStart + 1)
    $markRange.Delete()
}

# Two blank paragraphs after the existing trailing "AI" paragraph.
Add-EmptyPara $d
Add-EmptyPara $d

Add-TextPara $d "Att göra:"
Add-TextPara $d "Validera att alla enheter är utsatta när man trycker på start"
Add-TextPara $d "Fixa rörelselogiken"
Add-TextPara $d "Visa en vinst skärm när en flagga attackeras"
Add-TextPara $d "Dölj fieendepjäserna till man slåss"
Add-TextPara $d "AI"
Add-TextPara $d "EXTRA:"

# This paragraph is authored as three separate runs (the original edit split
# the word "HIGH" across two typing/save sessions). Build each chunk as its
# own paragraph, then merge them by removing the paragraph marks between
# them so the runs stay distinct instead of coalescing into one.
Add-TextPara $d "INZOOMNING PÅ ENHETERNA DÅ DE SLÅSS (RENDERA TVÅ STORA HI"
Add-TextPara $d "GH RES FILER"
Add-TextPara $d ")"
Merge-WithPrevious $d
Merge-WithPrevious $d

Add-TextPara $d "LJUD"

Write-Output "Final paragraph count: $($d.Paragraphs.Count)"
